$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8180572986602783
$ws.Range("B1").Value = 0.6014063358306885
$ws.Range("C1").Value = 0.4652432799339294
$ws.Range("D1").Value = 0.450622171163559
$ws.Range("E1").Value = 0.4864756464958191
